{"js": "// The tc/tcn/tl tag for this transcription's <id> was originally typed as\n// three separate runs: \"<id>\" (Courier New, brownish, 9pt), \"p161r_1\"\n// (plain black, default font), \"</id>\" (Courier New, brownish, 9pt).\n// Collapse them into a single run \"<id>p161r_1</id>\" that keeps the\n// formatting of the first run.\nconst body = context.document.body;\n\nconst results = body.search(\"<id>p161r_1</id>\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p161r_1</id>' text to merge.\");\n}\n\n// Replacing the whole (multi-run) match with identical text forces Word to\n// coalesce it into one run, carrying over the formatting of the first run\n// in the original range (Courier New / 7f6000 / 18 half-points).\nconst target = results.items[0];\ntarget.insertText(\"<id>p161r_1</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The tc/tcn/tl tag for this transcription's <id> was originally typed as\n# three separate runs: \"<id>\" (Courier New, brownish, 9pt), \"p161r_1\"\n# (plain black, default font), \"</id>\" (Courier New, brownish, 9pt).\n# Collapse them into a single run \"<id>p161r_1</id>\" that keeps the\n# formatting of the first run.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"<id>p161r_1</id>\")\n\nif (-not $found) {\n    throw \"Could not find '<id>p161r_1</id>' text to merge.\"\n}\n\n# Assigning text back onto a multi-run Range coalesces it into a single run\n# that inherits the formatting of the first run in the original range\n# (Courier New / 7f6000 / 18 half-points) \u2014 but only if the assigned text\n# actually differs from the current text, otherwise Word treats it as a\n# no-op and leaves the three runs untouched. Go through a scratch value\n# first so the final assignment is guaranteed to coalesce the runs.\n$rng.Text = \"<id>p161r_1</id>#\"\n$rng.Text = \"<id>p161r_1</id>\"\n"}
